# Swap the data of rows 10 and 11 for columns A, B, D, E, F, G, H, Q, R.
# (All other columns already hold identical values in both rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell10 = $ws.Range($col + "10")
    $cell11 = $ws.Range($col + "11")

    $val10 = $cell10.Value()
    $val11 = $cell11.Value()

    $cell10.Value = $val11
    $cell11.Value = $val10
}
